$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (rows 1-9 keep their row numbers) ---
$ws.Range("B3").Value = 0.0147
$ws.Range("A7").Value = "Iout ma"

# --- Insert a new row at position 10 ("Vc") ---
$ws.Rows("10").Insert()
$ws.Range("A10").Value = "Vc"
$ws.Range("B10").Formula = "=B9-B7*B4/1000"
$ws.Range("B10").Style = "Normal"

# Row 11 (old Rmax calc, previously row10) now needs to reference the new Vc cell (B10) instead of VDC (B9)
$ws.Range("B11").Formula = "=B10/B7*1000"

# The two pictures anchored below row 10 are free-floating and don't auto-shift with the row
# insert, so nudge them down by one default row height (15pt) to keep them visually aligned
# with the rows they used to sit next to.
$ws.Shapes.Item(3).Top = $ws.Shapes.Item(3).Top + 15
$ws.Shapes.Item(4).Top = $ws.Shapes.Item(4).Top + 15

# --- Update values further down (now shifted by +1 due to the row10 insert) ---
# Row 20 (old row19, R3 (Zin))
$ws.Range("B20").Value = 10000000

# Row 24 (old row23, R1 Actual)
$ws.Range("B24").Value = 950000

# --- New block rows 33-37 (labels typed in the order R1, R2, V1, Rdut to match shared-string order) ---
$ws.Range("D34").Value = "R1"
$ws.Range("E34").Value = 1000000
$ws.Range("E34").NumberFormat = "0.00E+00"

$ws.Range("D35").Value = "R2"
$ws.Range("E35").Value = 1400000
$ws.Range("E35").NumberFormat = "0.00E+00"

$ws.Range("D36").Value = "V1"
$ws.Range("E36").Formula = "=0.0147*E33*E35/(E34+E35)+B31"
$ws.Range("E36").NumberFormat = "0.00E+00"

$ws.Range("E37").NumberFormat = "0.00E+00"

$ws.Range("D33").Value = "Rdut"
$ws.Range("E33").Value = 100

# --- New row 31: Voff ---
$ws.Range("A31").Value = "Voff"
$ws.Range("B31").Value = 0.04

# --- Sheet view / selection tweaks ---
$ws.Range("K34").Select()
$excel.ActiveWindow.ScrollRow = 21

$excel.Calculate()
